# Insert a new column before column B (shifts existing B:P -> C:Q)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B").Insert()

# Fill the new column B (rows 2-17) with the L1/L2 list-membership labels
$values = @("L1","L2","L2","L1","L1","L2","L2","L1","L1","L2","L2","L1","L1","L2","L2","L1")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Header for the new column
$ws.Cells.Item(1, 2).Value = "list"
